# Auto-generated Excel COM-interop edit script
# Updates the cryptocurrency list prices and 1h volume/change values
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking numbers as TEXT in the source data
# (e.g. "141.08", thousand-dot formatted "57.317.20"). Force the whole column to
# Text format before writing so Excel does not auto-convert these into numbers,
# then restore the default "Normal" style so no stray formatting is left behind.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "57.317.20"
$ws.Range("E2").Value = "  +3.57%  "
$ws.Range("D3").Value = "3.064.91"
$ws.Range("E3").Value = "  +5.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "513.93"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "141.08"
$ws.Range("E6").Value = "  +6.66%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +6.93%  "
$ws.Range("D12").Value = "3.592.59"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "25.47"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("D16").Value = "57.366.51"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "3.069.64"
$ws.Range("E17").Value = "  +5.70%  "
$ws.Range("D18").Value = "5.98"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "8.13"
$ws.Range("E20").Value = "  +6.30%  "
$ws.Range("D21").Value = "337.13"
$ws.Range("E21").Value = "  +8.12%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "0.498"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "65.45"
$ws.Range("E24").Value = "  +4.81%  "
$ws.Range("E25").Value = "  +7.80%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +13.20%  "
$ws.Range("D28").Value = "6.41"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "7.06"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "20.73"
$ws.Range("E31").Value = "  +4.99%  "
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D33").Value = "154.39"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("D36").Value = "25.97"
$ws.Range("E36").Value = "  +7.80%  "
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("D38").Value = "0.0673"
$ws.Range("E38").Value = "  +4.88%  "
$ws.Range("D39").Value = "3.104.08"
$ws.Range("E39").Value = "  +5.69%  "
$ws.Range("D40").Value = "36.99"
$ws.Range("E40").Value = "  +2.28%  "

# Row 41 (was Mantle) becomes Filecoin; Row 42 (was Filecoin) becomes Mantle
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.669"
$ws.Range("E42").Value = "  +5.50%  "

$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "2.251.44"
$ws.Range("E44").Value = "  +7.86%  "
$ws.Range("E45").Value = "  +9.01%  "
$ws.Range("E46").Value = "  +5.13%  "
$ws.Range("D47").Value = "0.949"
$ws.Range("E47").Value = "  +4.90%  "
$ws.Range("E48").Value = "  +8.41%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "0.0866"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("E51").Value = "  +5.09%  "

# Restore default styling on the Price column now that the text values are set
$priceCol.Style = "Normal"

